# Generate Report for Handoff
#
# For the rows on the "zh-cn" / "de-de" localization-status sheets that are
# "Ready for handoff" and have not yet been picked up for handback (i.e. the
# handoff .xlf has just been (re-)generated), stamp the new handoff
# timestamp and mark the row's Priority as "ht" (handoff type).
#
# Rows affected (1-based sheet rows): 7, 8, 10, 12, 13, 14

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$rows = @(7, 8, 10, 12, 13, 14)

# New "Latest Handoff Datetime" stamps (handoff was regenerated a little
# later than the previous run captured in the workbook).
$newHandoffTimeOverviewDeDe = "2016-08-16 12:21:13"
$newHandoffTimeZhCn         = "2016-08-16 12:21:05"

foreach ($r in $rows) {
    # Overview sheet: column G = "Latest Handoff Datetime"
    $wsOverview.Range("G$r").Value = $newHandoffTimeOverviewDeDe

    # zh-cn sheet: column H = "Latest Handoff Datetime", column E = "Priority"
    $wsZhCn.Range("H$r").Value = $newHandoffTimeZhCn
    $wsZhCn.Range("E$r").Value = "ht"

    # de-de sheet: column H = "Latest Handoff Datetime", column E = "Priority"
    $wsDeDe.Range("H$r").Value = $newHandoffTimeOverviewDeDe
    $wsDeDe.Range("E$r").Value = "ht"
}
